# Lives Counter HUD update: re-highlight the "HUD" / "# of lives" /
# "Health bar" bullets (the ones currently marked yellow) in cyan.
$d = $word.ActiveDocument

$wdYellow = 7
$wdTurquoise = 3

$targets = @("HUD", "# of lives", "Health bar")

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()
    if (($targets -contains $text) -and ($p.Range.HighlightColorIndex -eq $wdYellow)) {
        $p.Range.Font.HighlightColorIndex = $wdTurquoise
    }
}

Write-Host "Re-highlighted HUD / # of lives / Health bar bullets to cyan"
